$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (rows 2-7) from 45183 to 45184
$ws.Range("C2:C7").Value = 45184
